# Refined metadata to be additional tab
# 1) Update "time_taken" (column F) timestamps on the "data" sheet to
#    reflect the new panel-query run time.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = "2021-10-05 14:33:18.022533"
$ws.Range("F3").Value = "2021-10-05 14:33:18.022542"
$ws.Range("F4").Value = "2021-10-05 14:33:18.022546"
$ws.Range("F5").Value = "2021-10-05 14:33:18.022549"
$ws.Range("F6").Value = "2021-10-05 14:33:18.022552"
$ws.Range("F7").Value = "2021-10-05 14:33:18.022555"
$ws.Range("F8").Value = "2021-10-05 14:33:18.022558"
$ws.Range("F9").Value = "2021-10-05 14:33:18.022560"
$ws.Range("F10").Value = "2021-10-05 14:33:18.022564"
$ws.Range("F11").Value = "2021-10-05 14:33:18.022567"
$ws.Range("F12").Value = "2021-10-05 14:33:18.022569"
$ws.Range("F13").Value = "2021-10-05 14:33:18.022572"
$ws.Range("F14").Value = "2021-10-05 14:33:18.022575"
$ws.Range("F15").Value = "2021-10-05 14:33:18.022578"
$ws.Range("F16").Value = "2021-10-05 14:33:18.022581"
$ws.Range("F17").Value = "2021-10-05 14:33:18.022583"
$ws.Range("F18").Value = "2021-10-05 14:33:18.022586"
$ws.Range("F19").Value = "2021-10-05 14:33:18.022589"
$ws.Range("F20").Value = "2021-10-05 14:33:18.022592"
$ws.Range("F21").Value = "2021-10-05 14:33:18.022594"
$ws.Range("F22").Value = "2021-10-05 14:33:18.022597"
$ws.Range("F23").Value = "2021-10-05 14:33:18.022600"
$ws.Range("F24").Value = "2021-10-05 14:33:18.022603"
$ws.Range("F25").Value = "2021-10-05 14:33:18.022605"
$ws.Range("F26").Value = "2021-10-05 14:33:18.022609"
$ws.Range("F27").Value = "2021-10-05 14:33:18.022611"
$ws.Range("F28").Value = "2021-10-05 14:33:18.022614"
$ws.Range("F29").Value = "2021-10-05 14:33:18.022617"
$ws.Range("F30").Value = "2021-10-05 14:33:18.022620"
$ws.Range("F31").Value = "2021-10-05 14:33:18.022622"
$ws.Range("F32").Value = "2021-10-05 14:33:18.022625"
$ws.Range("F33").Value = "2021-10-05 14:33:18.022627"
$ws.Range("F34").Value = "2021-10-05 14:33:18.022631"
$ws.Range("F35").Value = "2021-10-05 14:33:18.022634"
$ws.Range("F36").Value = "2021-10-05 14:33:18.022637"
$ws.Range("F37").Value = "2021-10-05 14:33:18.022639"
$ws.Range("F38").Value = "2021-10-05 14:33:18.022642"
$ws.Range("F39").Value = "2021-10-05 14:33:18.022645"

# 2) Add a new "metadata" worksheet, placed immediately after "data",
#    carrying the panel-query provenance info that used to live elsewhere.
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"
$meta.Range("A2").Value = 0

# Re-use the "data" sheet's header style (bold, centered, bordered) for the
# header row and the row-index cell, the same way the "data" sheet is styled.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B2").Value = "Brain Calcification"
$meta.Range("C2").Value = 58
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.12"
$meta.Range("E2").Value = "2021-08-23T03:16:23.882515Z"
$meta.Range("F2").Value = "2021-10-05 14:33:18.019159"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/58/?format=json"

# Keep "data" as the active/selected sheet (matches the original workbook's
# activeTab), since Worksheets.Add() would otherwise leave "metadata" active.
$ws.Activate()
